# Updates cryptos list (prices / 1h volume %) per GitHub Actions refresh.
# Values in column D that look numeric must be forced to Text format so
# Excel doesn't silently coerce them to numbers (which would also strip
# meaningful trailing zeros, e.g. "5.40" -> 5.4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Row 2 - Bitcoin
Set-TextCell "D2" "27.474.12"
$ws.Range("E2").Value = "  -0.86%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.617.71"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.12%  "

# Row 5 - BNB
Set-TextCell "D5" "211.01"
$ws.Range("E5").Value = "  -1.07%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -1.38%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.11%  "

# Row 8 - Solana
Set-TextCell "D8" "22.79"
$ws.Range("E8").Value = "  -1.13%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.262"
$ws.Range("E9").Value = "  +1.82%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.53%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextCell "D12" "1.847.43"
$ws.Range("E12").Value = "  -1.72%  "

# Row 13 - WrappedEther
Set-TextCell "D13" "1.624.25"
$ws.Range("E13").Value = "  -1.32%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.36%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -2.31%  "

# Row 16 - Litecoin
Set-TextCell "D16" "65.12"
$ws.Range("E16").Value = "  +1.38%  "

# Row 17 - WrappedBTC
Set-TextCell "D17" "27.468.85"
$ws.Range("E17").Value = "  -0.78%  "

# Row 18 - BitcoinCash
Set-TextCell "D18" "230.59"
$ws.Range("E18").Value = "  +0.15%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  -0.78%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -2.01%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.14%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.99%  "

# Row 23 - Avalanche
Set-TextCell "D23" "10.18"
$ws.Range("E23").Value = "  +1.61%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +5.91%  "

# Row 25 - Monero
Set-TextCell "D25" "150.48"
$ws.Range("E25").Value = "  +0.97%  "

# Row 26 / 27 - Cosmos and Stellar swapped rank order
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D26" "0.111"
$ws.Range("E26").Value = "  -1.18%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D27" "6.83"
$ws.Range("E27").Value = "  -1.93%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  +0.19%  "

# Row 29 - EthereumClassic
Set-TextCell "D29" "15.56"
$ws.Range("E29").Value = "  -0.65%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -1.00%  "

# Row 31 - Hedera
Set-TextCell "D31" "0.0482"
$ws.Range("E31").Value = "  -0.24%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -1.21%  "

# Row 33 - Maker
Set-TextCell "D33" "1.453.82"
$ws.Range("E33").Value = "  +0.90%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -3.37%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -3.49%  "

# Row 36 - HuobiToken
Set-TextCell "D36" "2.33"
$ws.Range("E36").Value = "  -0.27%  "

# Row 37 - TrustWalletToken
Set-TextCell "D37" "0.938"
$ws.Range("E37").Value = "  +3.82%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.42%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -2.12%  "

# Row 40 - ARBITRUM
Set-TextCell "D40" "0.862"
$ws.Range("E40").Value = "  -2.45%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.13%  "

# Row 42 - Aave
Set-TextCell "D42" "67.74"
$ws.Range("E42").Value = "  +3.26%  "

# Row 43 - mCoin
$ws.Range("E43").Value = "  +0.85%  "

# Row 44 - WEMIXToken
$ws.Range("E44").Value = "  -4.21%  "

# Row 45 - FraxShare
Set-TextCell "D45" "5.40"
$ws.Range("E45").Value = "  -4.65%  "

# Row 46 - MXToken
Set-TextCell "D46" "2.20"
$ws.Range("E46").Value = "  -2.52%  "

# Row 47 - RocketPoolETH
Set-TextCell "D47" "1.758.50"
$ws.Range("E47").Value = "  -1.70%  "

# Row 48 - RenderToken
Set-TextCell "D48" "1.71"
$ws.Range("E48").Value = "  +1.03%  "

# Row 49 - Quant
Set-TextCell "D49" "86.52"
$ws.Range("E49").Value = "  +0.21%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  +19.27%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  +1.72%  "
